$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows("13:13").Insert()

# The insert copies formatting from row 12 into the new row 13 col A; remove that stray cell
$ws.Range("A13").Clear()

# Populate new row 13 (B/C) with style copied from row 14 (B14 style=2, C14 style=3)
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value = '5840560 - Marco Antonio Carvalho Pereira'

# Fix text content that differs from the simple shift
$ws.Range("B10").Value = 'Fornecer uma visão geral da indústria de serviços. Desenvolver projeto de criação de um novo serviço.'
$ws.Range("C10").Value = 'Fornecer uma visão geral da indústria de serviços. Desenvolver projeto de criação de um novo serviço.'

$ws.Range("B14").Value = 'Introdução a Indústria de Serviços.Características Essenciais e diferenciadoras de Serviços.Projeto de Novo Serviço: Planejamento estratégico, Concepção do Serviço, Processos, Instalações. Avaliação e Melhoria.'
$ws.Range("C14").Value = 'Introdução a Indústria de Serviços.Características Essenciais e diferenciadoras de Serviços.Projeto de Novo Serviço: Planejamento estratégico, Concepção do Serviço, Processos, Instalações. Avaliação e Melhoria.'

$ws.Range("B16").Value = 'Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços.Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço(Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações deServiço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção eLocalização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria'
$ws.Range("C16").Value = 'Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços.Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço(Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações deServiço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção eLocalização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria'

$ws.Range("B19").Value = 'Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras.'
$ws.Range("C19").Value = 'Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras.'

$ws.Range("B20").Value = 'Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)'
$ws.Range("C20").Value = 'Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)'

$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação.'

$ws.Range("B22").Value = 'CARVALHO, M. M. (organizadora) e outros. Gestão de Serviços: Casos Brasileiros. Atlas. 2013CORREA, H. C. e CAON, M. Gestão de Serviços: Lucratividade por meio de operação e de satisfação dos clientes. Atlas, 2014FITZSIMMONS, J.; FITZSIMMONS, M.J. Administração de serviços: operações, estratégia e tecnologia de informação. PortoAlegre: Bookman, 2000.GIANESI, I e CORREA, H. Administração Estratégia de Serviços, ATLAS, 1995 – SPJOHNSTON, R. e CLARK, G. Administração e Operações de Serviços. Atlas, 2001LOVELOCK, C.H.; WRIGHT, L. Serviços: marketing e gestão. São Paulo: Saraiva, 2001.MELLO, C. H. P.; NETO, P. L. O. C.; TURRIONI, J.B. SILVA, C. E. S. Gestão do Processo de Desenvolvimento de Serviços.Atlas. 2010NORMANN, R. Administração de Serviços. São Paulo. Atlas. 1992.Bibliografia complementar será indicada ao longo do curso.'
$ws.Range("C22").Value = 'CARVALHO, M. M. (organizadora) e outros. Gestão de Serviços: Casos Brasileiros. Atlas. 2013CORREA, H. C. e CAON, M. Gestão de Serviços: Lucratividade por meio de operação e de satisfação dos clientes. Atlas, 2014FITZSIMMONS, J.; FITZSIMMONS, M.J. Administração de serviços: operações, estratégia e tecnologia de informação. PortoAlegre: Bookman, 2000.GIANESI, I e CORREA, H. Administração Estratégia de Serviços, ATLAS, 1995 – SPJOHNSTON, R. e CLARK, G. Administração e Operações de Serviços. Atlas, 2001LOVELOCK, C.H.; WRIGHT, L. Serviços: marketing e gestão. São Paulo: Saraiva, 2001.MELLO, C. H. P.; NETO, P. L. O. C.; TURRIONI, J.B. SILVA, C. E. S. Gestão do Processo de Desenvolvimento de Serviços.Atlas. 2010NORMANN, R. Administração de Serviços. São Paulo. Atlas. 1992.Bibliografia complementar será indicada ao longo do curso.'

# Column layout change: column A width now only spans col 1 (not 1:2)
$ws.Columns("1:1").ColumnWidth = 30.7109375

